{"js": "// The source text \"cunningly and with grace\" (spread across three runs)\n// becomes \"by cunning illusion\", keeping the formatting of the first run.\nconst body = context.document.body;\nconst results = body.search(\"cunningly and with grace\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase 'cunningly and with grace' not found\");\n}\n\nresults.items[0].insertText(\"by cunning illusion\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"cunningly and with grace\" (three runs) -> \"by cunning illusion\"\n# keeping the character formatting of the first run (\"cunningly\": rtl=0, no color).\n$d = $word.ActiveDocument\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"cunningly and with grace\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"by cunning illusion\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
